$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("error message")

# Insert a new row before row 8 (English block), shifting the existing
# English rows (old 8-13) down to 9-14. This makes room for the new
# "CH" (Chinese) row for beas_qc_1004 right after the other CH rows.
$ws.Rows("8:8").Insert()

# New CH row for beas_qc_1004 (row 8)
# Populate the new shared-string entries in the same order they appear
# in the target workbook: beas_qc_1004, then the English message, then
# the Chinese message.
$ws.Range("C8").Value = "beas_qc_1004"
$ws.Range("E15").Value = "Can't edit the qc position of this sample, please input the correct serial number."
$ws.Range("E8").Value = "不能输入该样品的质检项目，请输入正确的序列号"
$ws.Range("A8").Value = "CH  "
$ws.Range("B8").Value = "NULL"

# New E (English) row for beas_qc_1004 (row 15, appended at the end)
$ws.Range("A15").Value = "E   "
$ws.Range("B15").Value = "NULL"
$ws.Range("C15").Value = "beas_qc_1004"

# Column E width changed from autofit (bestFit) 100.375 to a fixed 92.
# (92 itself would round to 92.71 under the 7px Calibri-11 MDW rounding
# the host applies to ColumnWidth, so aim at the pre-rounding value that
# lands on exactly 92 once stored.)
$ws.Columns("E").ColumnWidth = 91.29

# Selection moved to E13.
$ws.Range("E13").Select()
